$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Walking Tests" header in column G, styled like the other header cells
# (wrap text, same as B1:F1).
$ws.Range("G1").Value = "Walking Tests"
$ws.Range("G1").WrapText = $true

# New walking-test score data for column G, rows 2-14.
$values = @(1, 1, 1, 1, 2, 2, 2, 2, 2, 1, 1, 1, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $values[$i]
}

# Widen the new column to fit its contents/header.
$ws.Columns.Item(7).ColumnWidth = 12.17

# Leave the selection where the author left it after adding the column.
$ws.Range("G13:H13").Select() | Out-Null
